$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 blank rows before the old "User knows how to navigate..." row (row 6),
# pushing it (and everything below) down to make room for the new header block.
$ws.Rows("6:9").Insert()

# The old assumption bullet "User knows how to navigate the website and add new
# veteranian" (now at row 10) is removed entirely - row becomes blank.
$ws.Range("A10").ClearContents()

# Update the remaining cell text content in place.
$ws.Range("A1").Value = "AC2: Security testing through SQL Injection"
$ws.Range("A4").Value = "SQL Injection tool (SQLMap) is installed"
$ws.Range("A5").Value = "Spring Pet Clinic website is connected to the server"
$ws.Range("A12").Value = "AC1TC12: SQL Injection on First Name input box"

# Update selection to match the author's saved cursor position.
$ws.Range("E4").Select()
